$d = $word.ActiveDocument

# ===========================================================================
# 1) Materials paragraph: "Rapamycin was purchased from XXX.  The fly stocks
#    were ordered from the Bloomington stock center (Table 1)" becomes
#    "Fly stocks were ordered from the Bloomington stock center or were
#    described previously (see Supplementary Table 1)"
# ===========================================================================
$d.Content.Find.Execute(
    "Rapamycin was purchased from XXX.  The fly stocks were ordered from the Bloomington stock center (Table 1)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Fly stocks were ordered from the Bloomington stock center or were described previously (see Supplementary Table 1)",
    2) | Out-Null

# ===========================================================================
# 2) Merge the Materials paragraph with the (empty) paragraph that used to
#    follow it, then append the new Rapamycin/Cayman sentence.
# ===========================================================================
$materialsIdx = -1
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text -like "Fly stocks were ordered from the Bloomington*") {
        $materialsIdx = $i
        break
    }
}

# Delete the paragraph mark ending the Materials paragraph -> merges it with
# the following empty paragraph (which disappears).
$mp = $d.Paragraphs($materialsIdx)
$d.Range($mp.Range.End - 1, $mp.Range.End).Delete() | Out-Null

# Append: "Rapamycin (Cayman Chemicals) was added where indicated when food
# was ~50<deg> C." with the degree sign in Cambria and the final "." bold.
$mp = $d.Paragraphs($materialsIdx)
$insPos = $mp.Range.End - 1
$d.Range($insPos, $insPos).InsertAfter("Rapamycin (Cayman Chemicals) was added where indicated when food was ~50") | Out-Null

$mp = $d.Paragraphs($materialsIdx)
$insPos = $mp.Range.End - 1
$d.Range($insPos, $insPos).InsertAfter([char]0x00B0) | Out-Null
$d.Range($insPos, $insPos + 1).Font.Name = "Cambria"

$mp = $d.Paragraphs($materialsIdx)
$insPos = $mp.Range.End - 1
$d.Range($insPos, $insPos).InsertAfter(" C") | Out-Null

$mp = $d.Paragraphs($materialsIdx)
$insPos = $mp.Range.End - 1
$d.Range($insPos, $insPos).InsertAfter(".") | Out-Null
$mp = $d.Paragraphs($materialsIdx)
$periodPos = $mp.Range.End - 1
$d.Range($periodPos - 1, $periodPos).Font.Bold = $true

# ===========================================================================
# 3) Restore the blank-paragraph count before "Tissue Culture and Myotube
#    Formation" (one spacer paragraph was consumed by the merge in step 2),
#    by inserting a new bold-formatted empty paragraph right before it.
# ===========================================================================
$tcIdx = -1
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text -like "Tissue Culture and Myotube Formation*") {
        $tcIdx = $i
        break
    }
}
$tp = $d.Paragraphs($tcIdx)
$hr = $tp.Range.Duplicate
$hr.Collapse(1)
$hr.InsertParagraphBefore() | Out-Null

$blankPara = $d.Paragraphs($tcIdx)
$blankPara.Range.Style = "Normal"
$blankPara.Range.Font.Bold = $true

# ===========================================================================
# 4) "Fly Breeding" section: the visible text is unchanged; only the
#    lastRenderedPageBreak rendering hint moved position (purely a cosmetic
#    layout marker with no text impact), so no action is required here.
# ===========================================================================

# ===========================================================================
# 5) Figure 5 legend rewording.
# ===========================================================================
$d.Content.Find.Execute(
    "Average climbing rate as measured during three intervals for",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Average climbing rate as measured during three age range intervals (in days) for",
    2) | Out-Null

$d.Content.Find.Execute(
    "relative to the control flies adjusted for multiple observations.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "relative to the control flies, and adjusted for multiple observations.",
    2) | Out-Null

# ===========================================================================
# 6) Append the new "Supplementary Tables and Figures" section right after
#    the "Figure 8 ... longevity." paragraph.
# ===========================================================================
$fig8Idx = -1
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text -like "Figure 8.*") {
        $fig8Idx = $i
    }
}

$fp = $d.Paragraphs($fig8Idx)
$r = $fp.Range.Duplicate
$r.Collapse(0)
$r.InsertParagraphAfter() | Out-Null

# 6a) Heading1 "Supplementary Tables and Figures"
$headingPara = $d.Paragraphs($fig8Idx + 1)
$headingPara.Range.Text = "Supplementary Tables and Figures"
$headingPara.Style = "Heading 1"

# 6b) A lone blank paragraph.
$r2 = $headingPara.Range.Duplicate
$r2.Collapse(0)
$r2.InsertParagraphAfter() | Out-Null
$blank2 = $d.Paragraphs($fig8Idx + 2)
$blank2.Style = "Normal"

# 6c) The Supplementary Figure 1 paragraph: bold title sentence followed by
#     normal body text.
$r3 = $blank2.Range.Duplicate
$r3.Collapse(0)
$r3.InsertParagraphAfter() | Out-Null
$sup1 = $d.Paragraphs($fig8Idx + 3)

$titleStart = $sup1.Range.Start
$d.Range($titleStart, $titleStart).InsertAfter("Supplementary Figure 1:  Dose response of rapamycin on fly eclosure and larvae development.") | Out-Null
$titleEnd = $d.Paragraphs($fig8Idx + 3).Range.End - 1
$d.Range($titleStart, $titleEnd).Font.Bold = $true

$bodyStart = $d.Paragraphs($fig8Idx + 3).Range.End - 1
$d.Range($bodyStart, $bodyStart).InsertAfter("  Flies were mated in the presence of varying doses of rapamycin.  After 7 days, the parental flies were removed.  After 21 days, the number of progeny and larval cases were counted.") | Out-Null
$bodyEnd = $d.Paragraphs($fig8Idx + 3).Range.End - 1
$d.Range($bodyStart, $bodyEnd).Font.Bold = $false

Write-Output "done"
